$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$usedRange.ClearContents()

$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"

$leftArr = New-Object 'object[,]' 7,8
$leftArr[0,0] = "name"
$leftArr[0,1] = "anchor score"
$leftArr[0,2] = "type occurences"
$leftArr[0,3] = "total occurences"
$leftArr[0,4] = "+%"
$leftArr[0,5] = "-%"
$leftArr[0,6] = "both"
$leftArr[0,7] = "normal"
$leftArr[1,0] = "crude"
$leftArr[1,1] = 0.7941176470588235
$leftArr[1,2] = 27
$leftArr[1,3] = 27
$leftArr[1,4] = 0
$leftArr[1,5] = 1
$leftArr[1,6] = $false
$leftArr[1,7] = 7
$leftArr[2,0] = "fraud"
$leftArr[2,1] = 0.7222222222222222
$leftArr[2,2] = 26
$leftArr[2,3] = 26
$leftArr[2,4] = 0
$leftArr[2,5] = 1
$leftArr[2,6] = $false
$leftArr[2,7] = 10
$leftArr[3,0] = "crisis"
$leftArr[3,1] = 0.5993150684931506
$leftArr[3,2] = 175
$leftArr[3,3] = 175
$leftArr[3,4] = 0
$leftArr[3,5] = 1
$leftArr[3,6] = $false
$leftArr[3,7] = 117
$leftArr[4,0] = "panic"
$leftArr[4,1] = 0.2151162790697674
$leftArr[4,2] = 111
$leftArr[4,3] = 111
$leftArr[4,4] = 0
$leftArr[4,5] = 1
$leftArr[4,6] = $false
$leftArr[4,7] = 405
$leftArr[5,0] = "sc"
$leftArr[5,1] = 0.1798941798941799
$leftArr[5,2] = 34
$leftArr[5,3] = 34
$leftArr[5,4] = 0
$leftArr[5,5] = 1
$leftArr[5,6] = $false
$leftArr[5,7] = 155
$leftArr[6,0] = "stop"
$leftArr[6,1] = 0.0992063492063492
$leftArr[6,2] = 25
$leftArr[6,3] = 25
$leftArr[6,4] = 0
$leftArr[6,5] = 1
$leftArr[6,6] = $false
$leftArr[6,7] = 227
$ws.Range("A2:H8").Value = $leftArr

$rightArr = New-Object 'object[,]' 43,8
$rightArr[0,0] = "name"
$rightArr[0,1] = "anchor score"
$rightArr[0,2] = "type occurences"
$rightArr[0,3] = "total occurences"
$rightArr[0,4] = "+%"
$rightArr[0,5] = "-%"
$rightArr[0,6] = "both"
$rightArr[0,7] = "normal"
$rightArr[1,0] = "strong"
$rightArr[1,1] = 0.9545454545454546
$rightArr[1,2] = 21
$rightArr[1,3] = 21
$rightArr[1,4] = 1
$rightArr[1,5] = 0
$rightArr[1,6] = $false
$rightArr[1,7] = 1
$rightArr[2,0] = "best"
$rightArr[2,1] = 0.9491525423728814
$rightArr[2,2] = 56
$rightArr[2,3] = 56
$rightArr[2,4] = 1
$rightArr[2,5] = 0
$rightArr[2,6] = $false
$rightArr[2,7] = 3
$rightArr[3,0] = "interesting"
$rightArr[3,1] = 0.9090909090909091
$rightArr[3,2] = 30
$rightArr[3,3] = 30
$rightArr[3,4] = 1
$rightArr[3,5] = 0
$rightArr[3,6] = $false
$rightArr[3,7] = 3
$rightArr[4,0] = "love"
$rightArr[4,1] = 0.8913043478260869
$rightArr[4,2] = 41
$rightArr[4,3] = 41
$rightArr[4,4] = 1
$rightArr[4,5] = 0
$rightArr[4,6] = $false
$rightArr[4,7] = 5
$rightArr[5,0] = "nice"
$rightArr[5,1] = 0.8518518518518519
$rightArr[5,2] = 23
$rightArr[5,3] = 23
$rightArr[5,4] = 1
$rightArr[5,5] = 0
$rightArr[5,6] = $false
$rightArr[5,7] = 4
$rightArr[6,0] = "happy"
$rightArr[6,1] = 0.8461538461538461
$rightArr[6,2] = 22
$rightArr[6,3] = 22
$rightArr[6,4] = 1
$rightArr[6,5] = 0
$rightArr[6,6] = $false
$rightArr[6,7] = 4
$rightArr[7,0] = "great"
$rightArr[7,1] = 0.8392857142857143
$rightArr[7,2] = 94
$rightArr[7,3] = 94
$rightArr[7,4] = 1
$rightArr[7,5] = 0
$rightArr[7,6] = $false
$rightArr[7,7] = 18
$rightArr[8,0] = "thanks"
$rightArr[8,1] = 0.8292682926829268
$rightArr[8,2] = 68
$rightArr[8,3] = 68
$rightArr[8,4] = 1
$rightArr[8,5] = 0
$rightArr[8,6] = $false
$rightArr[8,7] = 14
$rightArr[9,0] = "positive"
$rightArr[9,1] = 0.8103448275862069
$rightArr[9,2] = 47
$rightArr[9,3] = 47
$rightArr[9,4] = 1
$rightArr[9,5] = 0
$rightArr[9,6] = $false
$rightArr[9,7] = 11
$rightArr[10,0] = "free"
$rightArr[10,1] = 0.8
$rightArr[10,2] = 96
$rightArr[10,3] = 96
$rightArr[10,4] = 1
$rightArr[10,5] = 0
$rightArr[10,6] = $false
$rightArr[10,7] = 24
$rightArr[11,0] = "ensure"
$rightArr[11,1] = 0.7916666666666666
$rightArr[11,2] = 19
$rightArr[11,3] = 19
$rightArr[11,4] = 1
$rightArr[11,5] = 0
$rightArr[11,6] = $false
$rightArr[11,7] = 5
$rightArr[12,0] = "special"
$rightArr[12,1] = 0.7777777777777778
$rightArr[12,2] = 28
$rightArr[12,3] = 28
$rightArr[12,4] = 1
$rightArr[12,5] = 0
$rightArr[12,6] = $false
$rightArr[12,7] = 8
$rightArr[13,0] = "friends"
$rightArr[13,1] = 0.75
$rightArr[13,2] = 21
$rightArr[13,3] = 21
$rightArr[13,4] = 1
$rightArr[13,5] = 0
$rightArr[13,6] = $false
$rightArr[13,7] = 7
$rightArr[14,0] = "thank"
$rightArr[14,1] = 0.7421875
$rightArr[14,2] = 95
$rightArr[14,3] = 95
$rightArr[14,4] = 1
$rightArr[14,5] = 0
$rightArr[14,6] = $false
$rightArr[14,7] = 33
$rightArr[15,0] = "healthy"
$rightArr[15,1] = 0.7407407407407407
$rightArr[15,2] = 20
$rightArr[15,3] = 20
$rightArr[15,4] = 1
$rightArr[15,5] = 0
$rightArr[15,6] = $false
$rightArr[15,7] = 7
$rightArr[16,0] = "safe"
$rightArr[16,1] = 0.7183098591549296
$rightArr[16,2] = 102
$rightArr[16,3] = 102
$rightArr[16,4] = 1
$rightArr[16,5] = 0
$rightArr[16,6] = $false
$rightArr[16,7] = 40
$rightArr[17,0] = "support"
$rightArr[17,1] = 0.7075471698113207
$rightArr[17,2] = 75
$rightArr[17,3] = 75
$rightArr[17,4] = 1
$rightArr[17,5] = 0
$rightArr[17,6] = $false
$rightArr[17,7] = 31
$rightArr[18,0] = "safety"
$rightArr[18,1] = 0.7058823529411765
$rightArr[18,2] = 36
$rightArr[18,3] = 36
$rightArr[18,4] = 1
$rightArr[18,5] = 0
$rightArr[18,6] = $false
$rightArr[18,7] = 15
$rightArr[19,0] = "good"
$rightArr[19,1] = 0.65625
$rightArr[19,2] = 105
$rightArr[19,3] = 105
$rightArr[19,4] = 1
$rightArr[19,5] = 0
$rightArr[19,6] = $false
$rightArr[19,7] = 55
$rightArr[20,0] = "better"
$rightArr[20,1] = 0.6507936507936508
$rightArr[20,2] = 41
$rightArr[20,3] = 41
$rightArr[20,4] = 1
$rightArr[20,5] = 0
$rightArr[20,6] = $false
$rightArr[20,7] = 22
$rightArr[21,0] = "confidence"
$rightArr[21,1] = 0.6388888888888888
$rightArr[21,2] = 23
$rightArr[21,3] = 23
$rightArr[21,4] = 1
$rightArr[21,5] = 0
$rightArr[21,6] = $false
$rightArr[21,7] = 13
$rightArr[22,0] = "fresh"
$rightArr[22,1] = 0.6041666666666666
$rightArr[22,2] = 29
$rightArr[22,3] = 29
$rightArr[22,4] = 1
$rightArr[22,5] = 0
$rightArr[22,6] = $false
$rightArr[22,7] = 19
$rightArr[23,0] = "relief"
$rightArr[23,1] = 0.6
$rightArr[23,2] = 30
$rightArr[23,3] = 30
$rightArr[23,4] = 1
$rightArr[23,5] = 0
$rightArr[23,6] = $false
$rightArr[23,7] = 20
$rightArr[24,0] = "well"
$rightArr[24,1] = 0.5531914893617021
$rightArr[24,2] = 52
$rightArr[24,3] = 52
$rightArr[24,4] = 1
$rightArr[24,5] = 0
$rightArr[24,6] = $false
$rightArr[24,7] = 42
$rightArr[25,0] = "hand"
$rightArr[25,1] = 0.5143603133159269
$rightArr[25,2] = 197
$rightArr[25,3] = 197
$rightArr[25,4] = 1
$rightArr[25,5] = 0
$rightArr[25,6] = $false
$rightArr[25,7] = 186
$rightArr[26,0] = "heroes"
$rightArr[26,1] = 0.5106382978723404
$rightArr[26,2] = 24
$rightArr[26,3] = 24
$rightArr[26,4] = 1
$rightArr[26,5] = 0
$rightArr[26,6] = $false
$rightArr[26,7] = 23
$rightArr[27,0] = "like"
$rightArr[27,1] = 0.5058823529411764
$rightArr[27,2] = 172
$rightArr[27,3] = 172
$rightArr[27,4] = 1
$rightArr[27,5] = 0
$rightArr[27,6] = $false
$rightArr[27,7] = 168
$rightArr[28,0] = "care"
$rightArr[28,1] = 0.5056179775280899
$rightArr[28,2] = 45
$rightArr[28,3] = 45
$rightArr[28,4] = 1
$rightArr[28,5] = 0
$rightArr[28,6] = $false
$rightArr[28,7] = 44
$rightArr[29,0] = "important"
$rightArr[29,1] = 0.4444444444444444
$rightArr[29,2] = 20
$rightArr[29,3] = 20
$rightArr[29,4] = 1
$rightArr[29,5] = 0
$rightArr[29,6] = $false
$rightArr[29,7] = 25
$rightArr[30,0] = "help"
$rightArr[30,1] = 0.4440677966101695
$rightArr[30,2] = 131
$rightArr[30,3] = 131
$rightArr[30,4] = 1
$rightArr[30,5] = 0
$rightArr[30,6] = $false
$rightArr[30,7] = 164
$rightArr[31,0] = "protect"
$rightArr[31,1] = 0.410958904109589
$rightArr[31,2] = 30
$rightArr[31,3] = 30
$rightArr[31,4] = 1
$rightArr[31,5] = 0
$rightArr[31,6] = $false
$rightArr[31,7] = 43
$rightArr[32,0] = "hope"
$rightArr[32,1] = 0.4
$rightArr[32,2] = 26
$rightArr[32,3] = 26
$rightArr[32,4] = 1
$rightArr[32,5] = 0
$rightArr[32,6] = $false
$rightArr[32,7] = 39
$rightArr[33,0] = "increase"
$rightArr[33,1] = 0.3846153846153846
$rightArr[33,2] = 30
$rightArr[33,3] = 30
$rightArr[33,4] = 1
$rightArr[33,5] = 0
$rightArr[33,6] = $false
$rightArr[33,7] = 48
$rightArr[34,0] = "please"
$rightArr[34,1] = 0.3723849372384937
$rightArr[34,2] = 89
$rightArr[34,3] = 89
$rightArr[34,4] = 1
$rightArr[34,5] = 0
$rightArr[34,6] = $false
$rightArr[34,7] = 150
$rightArr[35,0] = "sure"
$rightArr[35,1] = 0.3125
$rightArr[35,2] = 20
$rightArr[35,3] = 20
$rightArr[35,4] = 1
$rightArr[35,5] = 0
$rightArr[35,6] = $false
$rightArr[35,7] = 44
$rightArr[36,0] = "store"
$rightArr[36,1] = 0.05257270693512305
$rightArr[36,2] = 47
$rightArr[36,3] = 47
$rightArr[36,4] = 1
$rightArr[36,5] = 0
$rightArr[36,6] = $false
$rightArr[36,7] = 847
$rightArr[37,0] = "online"
$rightArr[37,1] = 0.05023923444976076
$rightArr[37,2] = 21
$rightArr[37,3] = 21
$rightArr[37,4] = 1
$rightArr[37,5] = 0
$rightArr[37,6] = $false
$rightArr[37,7] = 397
$rightArr[38,0] = "shopping"
$rightArr[38,1] = 0.04567307692307692
$rightArr[38,2] = 19
$rightArr[38,3] = 19
$rightArr[38,4] = 1
$rightArr[38,5] = 0
$rightArr[38,6] = $false
$rightArr[38,7] = 397
$rightArr[39,0] = "grocery"
$rightArr[39,1] = 0.02996670366259711
$rightArr[39,2] = 27
$rightArr[39,3] = 27
$rightArr[39,4] = 1
$rightArr[39,5] = 0
$rightArr[39,6] = $false
$rightArr[39,7] = 874
$rightArr[40,0] = "19"
$rightArr[40,1] = 0.01308411214953271
$rightArr[40,2] = 28
$rightArr[40,3] = 31
$rightArr[40,4] = 0.9
$rightArr[40,5] = 0.09999999999999998
$rightArr[40,6] = $true
$rightArr[40,7] = 2112
$rightArr[41,0] = "co"
$rightArr[41,1] = 0.01128668171557562
$rightArr[41,2] = 35
$rightArr[41,3] = 41
$rightArr[41,4] = 0.85
$rightArr[41,5] = 0.15
$rightArr[41,6] = $true
$rightArr[41,7] = 3066
$rightArr[42,0] = "corona"
$rightArr[42,1] = 0.006572769953051643
$rightArr[42,2] = 21
$rightArr[42,3] = 27
$rightArr[42,4] = 0.78
$rightArr[42,5] = 0.22
$rightArr[42,6] = $true
$rightArr[42,7] = 3174
$ws.Range("J2:Q44").Value = $rightArr
